$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: columns L/M are formatted as Text (numFmt "@"). A plain .Value
# assignment on a Text-formatted cell is stored as a literal text string
# (matching real Excel typed-entry behaviour), but the source data in this
# sheet stores genuine numbers in those columns. Briefly switch the cell to
# a general number format, write the value, then restore the original
# "@" text format so the cell keeps its original style/appearance.
function Set-NumericValue {
    param($range, $value)
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = "@"
}

# Row 340 (2021-01-30): revised daily snapshot values
$ws.Range("E340").Value = 13
$ws.Range("F340").Value = 12
$ws.Range("G340").Value = 103

# Row 341 (2021-01-31)
$ws.Range("E341").Value = 14
$ws.Range("F341").Value = 11
$ws.Range("G341").Value = 112

# Row 342 (2021-02-01)
$ws.Range("E342").Value = 16
$ws.Range("F342").Value = 11
$ws.Range("G342").Value = 116

# Row 343 (2021-02-02)
$ws.Range("E343").Value = 15
$ws.Range("F343").Value = 9
$ws.Range("G343").Value = 121

# Row 344 (2021-02-03)
$ws.Range("C344").Value = 144
$ws.Range("E344").Value = 11
$ws.Range("F344").Value = 8
$ws.Range("G344").Value = 118

# Row 345 (2021-02-04)
$ws.Range("C345").Value = 93
$ws.Range("F345").Value = 8
$ws.Range("G345").Value = 114

# Row 346 (2021-02-05)
$ws.Range("C346").Value = 100
$ws.Range("F346").Value = 9
$ws.Range("G346").Value = 107
Set-NumericValue $ws.Range("L346") 2
Set-NumericValue $ws.Range("M346") 1

# Row 347 (2021-02-06) - newly reported data (previously empty)
$ws.Range("C347").Value = 47
$ws.Range("E347").Value = 12
$ws.Range("F347").Value = 9
$ws.Range("G347").Value = 107
Set-NumericValue $ws.Range("L347") 1
Set-NumericValue $ws.Range("M347") 0

# Row 348 (2021-02-07) - newly reported data (previously empty)
$ws.Range("C348").Value = 30
$ws.Range("E348").Value = 13
$ws.Range("F348").Value = 9
$ws.Range("G348").Value = 107
Set-NumericValue $ws.Range("L348") 1
Set-NumericValue $ws.Range("M348") 0

# Row 349 (2021-02-08) - newly reported data (previously empty)
$ws.Range("C349").Value = 13
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 9
$ws.Range("G349").Value = 108
Set-NumericValue $ws.Range("L349") 0
Set-NumericValue $ws.Range("M349") 0
